$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Map of (row, col) -> new text, matching the document's layout of
# 5 "problem" rows (1, 5, 9, 13, 17) each with 5 columns.
$updates = @(
    @{ Row = 1;  Col = 1; New = "94÷8=" },
    @{ Row = 1;  Col = 2; New = "54÷6=" },
    @{ Row = 1;  Col = 3; New = "93÷2=" },
    @{ Row = 1;  Col = 4; New = "77÷2=" },
    @{ Row = 1;  Col = 5; New = "83÷7=" },

    @{ Row = 5;  Col = 1; New = "85÷2=" },
    @{ Row = 5;  Col = 2; New = "76÷4=" },
    @{ Row = 5;  Col = 3; New = "44÷8=" },
    @{ Row = 5;  Col = 4; New = "70÷3=" },
    @{ Row = 5;  Col = 5; New = "95÷7=" },

    @{ Row = 9;  Col = 1; New = "59÷4=" },
    @{ Row = 9;  Col = 2; New = "93÷7=" },
    @{ Row = 9;  Col = 3; New = "43÷2=" },
    @{ Row = 9;  Col = 4; New = "66÷9=" },
    @{ Row = 9;  Col = 5; New = "99÷8=" },

    @{ Row = 13; Col = 1; New = "90÷2=" },
    @{ Row = 13; Col = 2; New = "15÷7=" },
    @{ Row = 13; Col = 3; New = "17÷4=" },
    @{ Row = 13; Col = 4; New = "11÷8=" },
    @{ Row = 13; Col = 5; New = "24÷7=" },

    @{ Row = 17; Col = 1; New = "57÷4=" },
    @{ Row = 17; Col = 2; New = "77÷2=" },
    @{ Row = 17; Col = 3; New = "84÷9=" },
    @{ Row = 17; Col = 4; New = "62÷7=" },
    @{ Row = 17; Col = 5; New = "22÷6=" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark/paragraph-mark characters so only the
    # visible text is replaced, preserving run formatting.
    $rng.End = $rng.End - 2
    $rng.Text = $u.New
}
